$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.703.44'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.914.89'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.43%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4919'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2995'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06764'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.917.41'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '17.16'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.60%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07326'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.209'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.51'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6739'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.686.66'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007966'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.53'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.99%  '
$ws.Range('E19').Value = '  +0.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.168.84'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.82%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.407'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +12.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.004'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '197.81'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.318'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.660'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '161.85'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.65'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.958'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.460'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.338'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09142'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.071'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05267'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7429'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.123'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.716'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01851'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.722'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9286'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.56%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.076'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.28%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4492'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.55%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '107.11'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.36%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '71.93'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +24.58%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.938'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.002'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1396'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.694'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.102'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.49'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05908'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4034'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.76%  '
